$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.582.85"
$ws.Range("E2").Value = "  -2.40%  "

$ws.Range("D3").Value = "1.818.22"
$ws.Range("E3").Value = "  -1.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.80%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4566"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.58%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3671"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07156"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8788"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07781"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.47%  "

$ws.Range("D13").Value = "1.776.98"
$ws.Range("E13").Value = "  -4.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.296"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.378"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.07%  "

$ws.Range("E17").Value = "  +0.82%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008614"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "

$ws.Range("D20").Value = "26.658.86"
$ws.Range("E20").Value = "  -2.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.68%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.42%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.990"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.084"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.871"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08699"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.066"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.531"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7366"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.714"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.38%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.74%  "

$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.084"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.93%  "

$ws.Range("E38").Value = "  +0.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05121"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.907"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.021"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5038"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.64%  "

$ws.Range("E43").Value = "  -4.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.200"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.007"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4631"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.985"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.57%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.595"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06027"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.50%  "
